# Applies the "cryptos list" refresh described in the commit diff:
# updates Price (D) / Volume(1h) (E) figures for the existing rows, and
# fixes rows 12/13 and 50/51 where the coin order (Name/Link/Price/Volume)
# had been swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.425.13"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.675.42"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'217.00"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").Value = "'0.5314"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("D8").Value = "'0.2694"
$ws.Range("E8").Value = "  +3.68%  "
$ws.Range("D9").Value = "'0.06404"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").Value = "'21.74"
$ws.Range("E10").Value = "  +4.97%  "
$ws.Range("D11").Value = "'0.07823"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.508"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.670.67"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").Value = "'0.5566"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "0.0₅8341"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").Value = "'65.64"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "26.462.26"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'4.730"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "'193.61"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").Value = "'6.340"
$ws.Range("E22").Value = "  +2.96%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'142.52"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").Value = "'0.1290"
$ws.Range("E25").Value = "  +5.96%  "
$ws.Range("D26").Value = "'7.395"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "'16.27"
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("D28").Value = "'1.439"
$ws.Range("E28").Value = "  +2.32%  "
$ws.Range("D29").Value = "'0.06358"
$ws.Range("E29").Value = "  +6.01%  "
$ws.Range("D30").Value = "'1.273"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "'3.630"
$ws.Range("E31").Value = "  +5.59%  "
$ws.Range("D32").Value = "'3.447"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").Value = "'0.6196"
$ws.Range("E35").Value = "  +8.16%  "
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("D37").Value = "'2.780"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("D38").Value = "'6.170"
$ws.Range("E38").Value = "  +7.69%  "
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").Value = "1.086.26"
$ws.Range("E40").Value = "  +4.81%  "
$ws.Range("D41").Value = "'0.8636"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D43").Value = "'100.31"
$ws.Range("D44").Value = "1.820.97"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").Value = "'57.37"
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("D46").Value = "'8.182"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").Value = "'0.05210"
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.484"
$ws.Range("E50").Value = "  +6.95%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'6.035"
$ws.Range("E51").Value = "  +1.90%  "
